$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("predictions")

# Row 3 (Group A, Italy vs Turkey) - final score 1-1 (draw, not yet reflected in Awon/Bwon/DrawHappened)
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1

# Row 6 (Group A, Switzerland vs Turkey) - draw happened, final score 1-1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 1

# Row 10 (Group B, Belgium vs Russia) - A won, final score 3-0
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("N10").Value = 3
$ws.Range("O10").Value = 0

# Row 11 (Group B, Denmark vs Finland) - B won, final score 0-1
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 1

# Update the active selection to match the saved view state
$ws.Range("J51").Select() | Out-Null
